# Apply the data change described by the diff:
# Column C ("Fitness") values are updated for rows 2-140 (data rows for Generation 0-138):
#   - Rows 2-10  (C2:C10)   : 7534 -> 7310
#   - Rows 11-140 (C11:C140): (7534/7345/7320/7295) -> 7293
# Rows 141 and below already contain 7293 and remain unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 through 10 become 7310
$ws.Range("C2:C10").Value = 7310

# Rows 11 through 140 become 7293
$ws.Range("C11:C140").Value = 7293
